# Scheduled runner update: refresh market-board derived profit figures
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1749.3334
$ws.Range("I40").Value = 1666.6666
$ws.Range("J40").Value = 1997.3334
$ws.Range("K40").Value = 1666.6666
$ws.Range("L40").Value = 1997.3334
$ws.Range("M40").Value = -1491.6666
$ws.Range("N40").Value = -2347.3334

$ws.Range("H62").Value = 3007.7222
$ws.Range("I62").Value = 1980
$ws.Range("J62").Value = 4292.375
$ws.Range("K62").Value = 1980
$ws.Range("L62").Value = 4292.375
$ws.Range("M62").Value = -1356
$ws.Range("N62").Value = -5540.375

$ws.Range("H65").Value = 3007.7222
$ws.Range("I65").Value = 1980
$ws.Range("J65").Value = 4292.375
$ws.Range("K65").Value = 9900
$ws.Range("L65").Value = 21461.875
$ws.Range("M65").Value = -6780
$ws.Range("N65").Value = -27701.875

$ws.Range("H70").Value = 51260.25
$ws.Range("I70").Value = 251075.5
$ws.Range("J70").Value = 1306.4375
$ws.Range("K70").Value = 753226.5
$ws.Range("L70").Value = 3919.3125
$ws.Range("M70").Value = -752956.5
$ws.Range("N70").Value = -4459.3125

$ws.Range("H73").Value = 51260.25
$ws.Range("I73").Value = 251075.5
$ws.Range("J73").Value = 1306.4375
$ws.Range("K73").Value = 753226.5
$ws.Range("L73").Value = 3919.3125
$ws.Range("M73").Value = -752290.5
$ws.Range("N73").Value = -5791.3125

$ws.Range("H98").Value = 3643.8718
$ws.Range("I98").Value = 2582.7334
$ws.Range("J98").Value = 7181
$ws.Range("K98").Value = 2582.7334
$ws.Range("L98").Value = 7181
$ws.Range("M98").Value = -1084.7334
$ws.Range("N98").Value = -10177

$ws.Range("H122").Value = 3643.8718
$ws.Range("I122").Value = 2582.7334
$ws.Range("J122").Value = 7181
$ws.Range("K122").Value = 7748.2002
$ws.Range("L122").Value = 21543
$ws.Range("M122").Value = -5298.2002
$ws.Range("N122").Value = -26443

$ws.Range("H139").Value = 60000
$ws.Range("J139").Value = 60000
$ws.Range("L139").Value = 60000
$ws.Range("N139").Value = -70280

$ws.Range("H140").Value = 79856
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 79856
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 79856
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -90216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10662.453
$ws.Range("I32").Value = 10766.456
$ws.Range("K32").Value = 10766.456
$ws.Range("M32").Value = -10479.456

$ws.Range("H122").Value = 1734
$ws.Range("I122").Value = 2101
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 6303
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -3853
$ws.Range("N122").Value = -7900

$ws.Range("H131").Value = 55235.668
$ws.Range("J131").Value = 55235.668
$ws.Range("L131").Value = 55235.668
$ws.Range("N131").Value = -65315.668

$ws.Range("H132").Value = 3827.5293
$ws.Range("I132").Value = 4671.919
$ws.Range("K132").Value = 14015.757
$ws.Range("M132").Value = -11485.757

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5040.933
$ws.Range("I134").Value = 1724.6383
$ws.Range("J134").Value = 17030.615
$ws.Range("K134").Value = 5173.9149
$ws.Range("L134").Value = 51091.845
$ws.Range("M134").Value = -2638.9149
$ws.Range("N134").Value = -56161.845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 4252.1665
$ws.Range("I21").Value = 2013
$ws.Range("J21").Value = 4700
$ws.Range("K21").Value = 2013
$ws.Range("L21").Value = 4700
$ws.Range("M21").Value = -1778
$ws.Range("N21").Value = -5170

$ws.Range("H131").Value = 24500
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 24500
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 24500
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -34580

$ws.Range("H132").Value = 209006.2
$ws.Range("I132").Value = 265785.9
$ws.Range("J132").Value = 2165.7856
$ws.Range("K132").Value = 797357.7000000001
$ws.Range("L132").Value = 6497.3568
$ws.Range("M132").Value = -794827.7000000001
$ws.Range("N132").Value = -11557.3568

$ws.Range("H134").Value = 1392.0435
$ws.Range("I134").Value = 1091.5151
$ws.Range("J134").Value = 2154.923
$ws.Range("K134").Value = 3274.5453
$ws.Range("L134").Value = 6464.768999999999
$ws.Range("M134").Value = -739.5453000000002
$ws.Range("N134").Value = -11534.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 8145.3687
$ws.Range("I87").Value = 1085.6666
$ws.Range("J87").Value = 9469.0625
$ws.Range("K87").Value = 3256.9998
$ws.Range("L87").Value = 28407.1875
$ws.Range("M87").Value = -2008.9998
$ws.Range("N87").Value = -30903.1875

$ws.Range("H90").Value = 8145.3687
$ws.Range("I90").Value = 1085.6666
$ws.Range("J90").Value = 9469.0625
$ws.Range("K90").Value = 9770.999400000001
$ws.Range("L90").Value = 85221.5625
$ws.Range("M90").Value = -3530.999400000001
$ws.Range("N90").Value = -97701.5625

$ws.Range("H92").Value = 479.4
$ws.Range("I92").Value = 504.66666
$ws.Range("J92").Value = 441.5
$ws.Range("K92").Value = 1513.99998
$ws.Range("L92").Value = 1324.5
$ws.Range("M92").Value = -265.9999800000001
$ws.Range("N92").Value = -3820.5

$ws.Range("H98").Value = 1619.7
$ws.Range("I98").Value = 2977.25
$ws.Range("J98").Value = 714.6667
$ws.Range("K98").Value = 8931.75
$ws.Range("L98").Value = 2144.0001
$ws.Range("M98").Value = -7433.75
$ws.Range("N98").Value = -5140.0001

$ws.Range("H131").Value = 5303.1875
$ws.Range("I131").Value = 11513.223
$ws.Range("J131").Value = 2873.1738
$ws.Range("K131").Value = 34539.669
$ws.Range("L131").Value = 8619.5214
$ws.Range("M131").Value = -29499.669
$ws.Range("N131").Value = -18699.5214

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 13515.125
$ws.Range("J51").Value = 13515.125
$ws.Range("L51").Value = 13515.125
$ws.Range("N51").Value = -14533.125

$ws.Range("H109").Value = 13451.25
$ws.Range("J109").Value = 13451.25
$ws.Range("L109").Value = 13451.25
$ws.Range("N109").Value = -15531.25

$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H123").Value = 9817.684999999999
$ws.Range("J123").Value = 9817.684999999999
$ws.Range("L123").Value = 9817.684999999999
$ws.Range("N123").Value = -14717.685

$ws.Range("H131").Value = 31349
$ws.Range("J131").Value = 31349
$ws.Range("L131").Value = 31349
$ws.Range("N131").Value = -41429

$ws.Range("H132").Value = 3285.3062
$ws.Range("I132").Value = 3134.525
$ws.Range("J132").Value = 3955.4443
$ws.Range("K132").Value = 9403.575000000001
$ws.Range("L132").Value = 11866.3329
$ws.Range("M132").Value = -6873.575000000001
$ws.Range("N132").Value = -16926.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3715.276
$ws.Range("I132").Value = 3140
$ws.Range("J132").Value = 5225.375
$ws.Range("K132").Value = 9420
$ws.Range("L132").Value = 15676.125
$ws.Range("M132").Value = -6890
$ws.Range("N132").Value = -20736.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 24765.479
$ws.Range("J123").Value = 24765.479
$ws.Range("L123").Value = 24765.479
$ws.Range("N123").Value = -34565.479

$ws.Range("H125").Value = 53332.668
$ws.Range("J125").Value = 53332.668
$ws.Range("L125").Value = 53332.668
$ws.Range("N125").Value = -63172.668
